$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_4_4_0"
$ws.Range("B2").Value = 0.6684557879903787
$ws.Range("C2").Value = 0.9100156218193028
$ws.Range("D2").Value = 0.2829972253910134
$ws.Range("E2").Value = 0.509906072203437
$ws.Range("F2").Value = 0.3669217526912689
$ws.Range("G2").Value = 0.1010679975152016
$ws.Range("H2").Value = 1.664248585700989
$ws.Range("I2").Value = 0.8366822600364685

$ws.Range("A3").Value = "model_4_4_1"
$ws.Range("B3").Value = 0.683626272452579
$ws.Range("C3").Value = 0.9049237008223188
$ws.Range("D3").Value = 0.3042575693184714
$ws.Range("E3").Value = 0.521735952124135
$ws.Range("F3").Value = 0.3501324653625488
$ws.Range("G3").Value = 0.1067871153354645
$ws.Range("H3").Value = 1.614900827407837
$ws.Range("I3").Value = 0.8164864182472229

$ws.Range("A4").Value = "model_4_4_2"
$ws.Range("B4").Value = 0.6975153306565182
$ws.Range("C4").Value = 0.8947925582086143
$ws.Range("D4").Value = 0.3221677590490467
$ws.Range("E4").Value = 0.5296655678447494
$ws.Range("F4").Value = 0.334761381149292
$ws.Range("G4").Value = 0.1181661263108253
$ws.Range("H4").Value = 1.573329210281372
$ws.Range("I4").Value = 0.802949070930481

$ws.Range("A5").Value = "model_4_4_11"
$ws.Range("B5").Value = 0.6994980780308234
$ws.Range("C5").Value = 0.8742659282225739
$ws.Range("D5").Value = 0.3229044523762951
$ws.Range("E5").Value = 0.5229877998731991
$ws.Range("F5").Value = 0.3325670659542084
$ws.Range("G5").Value = 0.1412210911512375
$ws.Range("H5").Value = 1.571619153022766
$ws.Range("I5").Value = 0.8143492937088013

$ws.Range("A6").Value = "model_4_4_13"
$ws.Range("B6").Value = 0.7020053903228777
$ws.Range("C6").Value = 0.8728205407014865
$ws.Range("D6").Value = 0.3303691317837597
$ws.Range("E6").Value = 0.5272603113083143
$ws.Range("F6").Value = 0.3297922313213348
$ws.Range("G6").Value = 0.1428444981575012
$ws.Range("H6").Value = 1.554292678833008
$ws.Range("I6").Value = 0.8070554137229919

$ws.Range("A7").Value = "model_4_4_10"
$ws.Range("B7").Value = 0.7022782646560677
$ws.Range("C7").Value = 0.8761049879829474
$ws.Range("D7").Value = 0.3326048649968711
$ws.Range("E7").Value = 0.5298349473314377
$ws.Range("F7").Value = 0.3294902145862579
$ws.Range("G7").Value = 0.1391555070877075
$ws.Range("H7").Value = 1.549103379249573
$ws.Range("I7").Value = 0.8026600480079651

$ws.Range("A8").Value = "model_4_4_9"
$ws.Range("B8").Value = 0.7026112738418333
$ws.Range("C8").Value = 0.8764254340771886
$ws.Range("D8").Value = 0.3339027982645327
$ws.Range("E8").Value = 0.5307767252835462
$ws.Range("F8").Value = 0.3291216492652893
$ws.Range("G8").Value = 0.1387955844402313
$ws.Range("H8").Value = 1.546090722084045
$ws.Range("I8").Value = 0.8010520935058594

$ws.Range("A9").Value = "model_4_4_12"
$ws.Range("B9").Value = 0.7027623781178884
$ws.Range("C9").Value = 0.8737188735154413
$ws.Range("D9").Value = 0.3324128347244635
$ws.Range("E9").Value = 0.5288803183896613
$ws.Range("F9").Value = 0.3289544880390167
$ws.Range("G9").Value = 0.1418355256319046
$ws.Range("H9").Value = 1.549549102783203
$ws.Range("I9").Value = 0.8042897582054138

$ws.Range("A10").Value = "model_4_4_24"
$ws.Range("B10").Value = 0.7029256953378917
$ws.Range("C10").Value = 0.8780341444304357
$ws.Range("D10").Value = 0.3288350033750518
$ws.Range("E10").Value = 0.5280949423885281
$ws.Range("F10").Value = 0.3287737369537354
$ws.Range("G10").Value = 0.1369887292385101
$ws.Range("H10").Value = 1.557853698730469
$ws.Range("I10").Value = 0.8056304454803467

$ws.Range("A11").Value = "model_4_4_23"
$ws.Range("B11").Value = 0.7029373400039836
$ws.Range("C11").Value = 0.8772191374683332
$ws.Range("D11").Value = 0.3294285461558285
$ws.Range("E11").Value = 0.5281905321196545
$ws.Range("F11").Value = 0.3287608325481415
$ws.Range("G11").Value = 0.1379041224718094
$ws.Range("H11").Value = 1.55647599697113
$ws.Range("I11").Value = 0.8054672479629517

$ws.Range("A12").Value = "model_4_4_19"
$ws.Range("B12").Value = 0.7030944778229764
$ws.Range("C12").Value = 0.8774197168237126
$ws.Range("D12").Value = 0.3300992948189122
$ws.Range("E12").Value = 0.5286901239716139
$ws.Range("F12").Value = 0.3285869061946869
$ws.Range("G12").Value = 0.1376788169145584
$ws.Range("H12").Value = 1.554919123649597
$ws.Range("I12").Value = 0.8046144843101501

$ws.Range("A13").Value = "model_4_4_21"
$ws.Range("B13").Value = 0.7031227151132607
$ws.Range("C13").Value = 0.8773821715584417
$ws.Range("D13").Value = 0.3300855849611407
$ws.Range("E13").Value = 0.5286679213634493
$ws.Range("F13").Value = 0.3285556733608246
$ws.Range("G13").Value = 0.1377210021018982
$ws.Range("H13").Value = 1.554951071739197
$ws.Range("I13").Value = 0.8046522736549377

$ws.Range("A14").Value = "model_4_4_22"
$ws.Range("B14").Value = 0.7031408264771586
$ws.Range("C14").Value = 0.8772703832398968
$ws.Range("D14").Value = 0.3301359404447345
$ws.Range("E14").Value = 0.5286608344807811
$ws.Range("F14").Value = 0.3285356163978577
$ws.Range("G14").Value = 0.1378465592861176
$ws.Range("H14").Value = 1.554834008216858
$ws.Range("I14").Value = 0.8046643733978271

$ws.Range("A15").Value = "model_4_4_18"
$ws.Range("B15").Value = 0.7031994715781618
$ws.Range("C15").Value = 0.8776579788133476
$ws.Range("D15").Value = 0.3304066478604741
$ws.Range("E15").Value = 0.5289692955899117
$ws.Range("F15").Value = 0.3284707069396973
$ws.Range("G15").Value = 0.1374112367630005
$ws.Range("H15").Value = 1.554205656051636
$ws.Range("I15").Value = 0.8041377663612366

$ws.Range("A16").Value = "model_4_4_20"
$ws.Range("B16").Value = 0.70322471651263
$ws.Range("C16").Value = 0.8773515443555544
$ws.Range("D16").Value = 0.3305039993319383
$ws.Range("E16").Value = 0.5289248437376517
$ws.Range("F16").Value = 0.32844278216362
$ws.Range("G16").Value = 0.137755423784256
$ws.Range("H16").Value = 1.553979754447937
$ws.Range("I16").Value = 0.8042136430740356

$ws.Range("A17").Value = "model_4_4_14"
$ws.Range("B17").Value = 0.7037086977801366
$ws.Range("C17").Value = 0.8758545126375682
$ws.Range("D17").Value = 0.3338394054773862
$ws.Range("E17").Value = 0.5305368056836393
$ws.Range("F17").Value = 0.3279071748256683
$ws.Range("G17").Value = 0.1394368410110474
$ws.Range("H17").Value = 1.546237826347351
$ws.Range("I17").Value = 0.8014616966247559

$ws.Range("A18").Value = "model_4_4_17"
$ws.Range("B18").Value = 0.7037204078901853
$ws.Range("C18").Value = 0.8777074779903243
$ws.Range("D18").Value = 0.3323417524304575
$ws.Range("E18").Value = 0.530224291075495
$ws.Range("F18").Value = 0.3278942108154297
$ws.Range("G18").Value = 0.1373556405305862
$ws.Range("H18").Value = 1.549714088439941
$ws.Range("I18").Value = 0.8019952774047852

$ws.Range("A19").Value = "model_4_4_15"
$ws.Range("B19").Value = 0.7040584575013296
$ws.Range("C19").Value = 0.8772512338560065
$ws.Range("D19").Value = 0.3339609061715154
$ws.Range("E19").Value = 0.5311021211286993
$ws.Range("F19").Value = 0.3275201022624969
$ws.Range("G19").Value = 0.1378680765628815
$ws.Range("H19").Value = 1.545955896377563
$ws.Range("I19").Value = 0.8004966974258423

$ws.Range("A20").Value = "model_4_4_16"
$ws.Range("B20").Value = 0.7042708279821837
$ws.Range("C20").Value = 0.8777545766189556
$ws.Range("D20").Value = 0.3342810293407555
$ws.Range("E20").Value = 0.5314817303507453
$ws.Range("F20").Value = 0.3272850811481476
$ws.Range("G20").Value = 0.1373027414083481
$ws.Range("H20").Value = 1.545212864875793
$ws.Range("I20").Value = 0.7998486161231995

$ws.Range("A21").Value = "model_4_4_3"
$ws.Range("B21").Value = 0.7061710945194537
$ws.Range("C21").Value = 0.8932996562493819
$ws.Range("D21").Value = 0.3461915087233728
$ws.Range("E21").Value = 0.5445166343783956
$ws.Range("F21").Value = 0.3251820206642151
$ws.Range("G21").Value = 0.1198429241776466
$ws.Range("H21").Value = 1.517567157745361
$ws.Range("I21").Value = 0.7775956392288208

$ws.Range("A22").Value = "model_4_4_8"
$ws.Range("B22").Value = 0.7062253158344622
$ws.Range("C22").Value = 0.8815263157746647
$ws.Range("D22").Value = 0.3435388818733195
$ws.Range("E22").Value = 0.5387185307845261
$ws.Range("F22").Value = 0.3251220285892487
$ws.Range("G22").Value = 0.1330664157867432
$ws.Range("H22").Value = 1.523724317550659
$ws.Range("I22").Value = 0.7874940633773804

$ws.Range("A23").Value = "model_4_4_6"
$ws.Range("B23").Value = 0.708713968026623
$ws.Range("C23").Value = 0.8856122377451789
$ws.Range("D23").Value = 0.351049716694258
$ws.Range("E23").Value = 0.5449473662129757
$ws.Range("F23").Value = 0.3223678171634674
$ws.Range("G23").Value = 0.1284772157669067
$ws.Range("H23").Value = 1.506290674209595
$ws.Range("I23").Value = 0.776860237121582

$ws.Range("A24").Value = "model_4_4_4"
$ws.Range("B24").Value = 0.7094371083248838
$ws.Range("C24").Value = 0.8891990247186087
$ws.Range("D24").Value = 0.3574031548441156
$ws.Range("E24").Value = 0.5502623409308316
$ws.Range("F24").Value = 0.3215675055980682
$ws.Range("G24").Value = 0.1244486272335052
$ws.Range("H24").Value = 1.491543650627136
$ws.Range("I24").Value = 0.7677865028381348

$ws.Range("A25").Value = "model_4_4_7"
$ws.Range("B25").Value = 0.7114484381446619
$ws.Range("C25").Value = 0.8846762934229619
$ws.Range("D25").Value = 0.359859779531388
$ws.Range("E25").Value = 0.5502582858516588
$ws.Range("F25").Value = 0.3193415999412537
$ws.Range("G25").Value = 0.1295284628868103
$ws.Range("H25").Value = 1.485841512680054
$ws.Range("I25").Value = 0.7677934765815735

$ws.Range("A26").Value = "model_4_4_5"
$ws.Range("B26").Value = 0.7164766403374467
$ws.Range("C26").Value = 0.8927093415060943
$ws.Range("D26").Value = 0.3768346786395086
$ws.Range("E26").Value = 0.563916923414316
$ws.Range("F26").Value = 0.3137767910957336
$ws.Range("G26").Value = 0.1205059587955475
$ws.Range("H26").Value = 1.446440815925598
$ws.Range("I26").Value = 0.7444756031036377
